$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "ICe Cream"
$ws.Range("C3").Value = 200

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "straw"
$ws.Range("C4").Value = 32
